# Mise a jour de la navigation
# Update existing row 21 (EntrepriseController / ListeEntreprise search entry) and
# add two new navigation rows (22: liste entreprises pour modification, 23: autorisation d'absence)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: complete the existing "Recherche des entreprises" entry -------------
$ws.Range("B21").Value = "Recherche des entreprises qui ont déjà pris des stagiaires"
$ws.Range("C21").Value = "EntrepriseController"
$ws.Range("D21").Value = "ListeEntreprise "

# New font (size 12 Calibri) used for the Controller/Action columns of the new rows
$ws.Range("C21").Font.Size = 12
$ws.Range("D21").Font.Size = 12
$ws.Rows.Item(21).RowHeight = 15.75

# --- Row 22: new "Liste des entreprises pour modification" entry -----------------
$ws.Range("A22").Value = "Collaborateur"
$ws.Range("B22").Value = "Liste des entreprises pour modification/creation,suppression"
$ws.Range("C22").Value = "EntrepriseController"
$ws.Range("D22").Value = "ListeEntreprisePourModification"

$ws.Range("C22").Font.Size = 12
$ws.Range("D22").Font.Size = 12
$ws.Rows.Item(22).RowHeight = 15.75

# --- Row 23: new "Saisie autorisation d'absence" entry ---------------------------
$ws.Range("A23").Value = "Sta"
$ws.Range("B23").Value = "Saisie autorisation d'absence"
$ws.Range("C23").Value = "AutorisationAbsenceController"
$ws.Range("D23").Value = "CompleterInfoAbsence "

$ws.Range("C23").Font.Size = 12
$ws.Range("D23").Font.Size = 12
$ws.Rows.Item(23).RowHeight = 15.75

# --- Column widths: new column A, widened column B -------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.140625
$ws.Columns.Item(2).ColumnWidth = 56.140625

# --- Selection follows the last edited cell --------------------------------------
$ws.Range("D23").Select() | Out-Null
